$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.358.61"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "2.410.21"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'561.04"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "'135.83"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "'24.76"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "2.840.84"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "60.262.39"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "2.416.23"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "'325.33"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'64.60"
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("D24").Value = "'0.175"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "'8.50"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "'1.82"
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "'170.87"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +6.53%  "
$ws.Range("D33").Value = "'0.404"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "'18.35"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.60"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'323.46"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("D41").Value = "'38.79"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("D42").Value = "'147.33"
$ws.Range("E42").Value = "  +5.72%  "
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").Value = "'0.0970"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "'19.93"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D47").Value = "'0.574"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").Value = "'11.06"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  -0.54%  "
